$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) "Toggles del" sheet: update existing toggle counts and append a 4th row
# ---------------------------------------------------------------------------
$wsDel = $wb.Worksheets.Item("Toggles del")
$wsDel.Range("B2").Value = 72
$wsDel.Range("B3").Value = 120
$wsDel.Range("B4").Value = 56

# New row 5 (index 3) - copy the formatting used by the existing index rows
# (A2:A4) then set the values for A5/B5.
$wsDel.Range("A5").Value = 3
$wsDel.Range("B5").Value = 8
$wsDel.Range("A4").Copy()
$wsDel.Range("A5").PasteSpecial(-4122)

# Point the existing bar chart's series at the now-larger range so the
# chart formula covers B2:B5 instead of B2:B4.
$chartDel = $wsDel.ChartObjects(1).Chart
$chartDel.SeriesCollection(1).Values = "='Toggles del'!`$B`$2:`$B`$5"

# ---------------------------------------------------------------------------
# 2) "Toggles input del" sheet: identical change as "Toggles del"
# ---------------------------------------------------------------------------
$wsInputDel = $wb.Worksheets.Item("Toggles input del")
$wsInputDel.Range("B2").Value = 72
$wsInputDel.Range("B3").Value = 120
$wsInputDel.Range("B4").Value = 56

$wsInputDel.Range("A5").Value = 3
$wsInputDel.Range("B5").Value = 8
$wsInputDel.Range("A4").Copy()
$wsInputDel.Range("A5").PasteSpecial(-4122)

$chartInputDel = $wsInputDel.ChartObjects(1).Chart
$chartInputDel.SeriesCollection(1).Values = "='Toggles input del'!`$B`$2:`$B`$5"

# ---------------------------------------------------------------------------
# 3) "Correlation matrix" sheet: recomputed correlation values for rows 3/4
#    (both rows share identical new values)
# ---------------------------------------------------------------------------
$wsCorr = $wb.Worksheets.Item("Correlation matrix")
foreach ($r in 3, 4) {
    $wsCorr.Range("B$r").Value = 0.3162277660168379
    $wsCorr.Range("C$r").Value = 0
    $wsCorr.Range("D$r").Value = 0.1825741858350554
    $wsCorr.Range("E$r").Value = 0.3162277660168379
    $wsCorr.Range("F$r").Value = 0.3162277660168379
    $wsCorr.Range("G$r").Value = 0.4898979485566357
}
